# A new review row was inserted above the existing data (the old row 2
# "com.hamxa.shaynachim / bitcoin / cohenyossi408@gmail.com / cohenn167@gmail.com /
#  27/5/2019 15:59 / awesome app with great addictive concept / no" and old row 3
# both shift down by one, and the trailing blank row shifts from row 5 to row 6).
#
# Insert a brand new row at row 2 (this pushes everything below down one row,
# dimension grows from A1:G5 to A1:G6 automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(2).Insert()

# Give the new row 2 the same formatting (styles/number formats/height) as
# the data row right below it (old row 2, now row 3) before filling in values.
$ws.Range("A3:G3").Copy()
$ws.Range("A2:G2").PasteSpecial(-4122)
$ws.Rows(2).RowHeight = 13.8

# Fill in the new review's data.
$ws.Range("A2").Value = "com.hamxa.shaynachim"
$ws.Range("B2").Value = "bitcoin"
$ws.Range("C2").Value = "itamaramir2@gmail.com"
$ws.Range("D2").Value = "cohenyossi408@gmail.com"
$ws.Range("E2").Value = "27/5/2019 15:59"
$ws.Range("F2").Value = "sweet car albama! Hahaha"
$ws.Range("G2").Value = "no"

# Leave the selection where the user last typed (the new review text cell),
# matching the saved view state (activeCell F2, sqref F2).
[void]$ws.Range("F2").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
